$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.687.56"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.896.61"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.40%  "
$ws.Range("D7").Value = "3.897.01"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("E13").Value = "  +4.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").Value = "4.547.38"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "3.881.38"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "69.671.17"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").Value = "4.043.81"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "3.861.39"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "434.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("E48").Value = "  +20.93%  "
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
